# 6 hours by turn fix
# The schedule sheet needs its "Almoço" (lunch) break moved one slot earlier
# (12:20 instead of 13:00), which shifts every subsequent time label and
# activity down by one row and pushes the final rows of the table down,
# extending the grid from A1:F14 to A1:F17. Also fixes a swapped pair of
# class entries in C3/C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the swapped "Pedro-Resistencia mecanica" entry between C3 and C4 ---
$ws.Range("C3").Value = "Pedro-Resistencia mecanica"
$ws.Range("C4").Value = "-"

# --- Insert a new blank row above row 9 so the whole lower block shifts down ---
# This pushes old row 9 (13:00 row) down to row 10, old row 10 down to row 11,
# etc, and creates space for the new 12:20 "Almoço" row.
$ws.Rows("9:9").Insert()

# --- Row 8: lunch ("Almoço") moves out of this row; becomes a normal "-" slot ---
$ws.Cells.Item(8, 2).Value = "-"
$ws.Cells.Item(8, 3).Value = "-"
$ws.Cells.Item(8, 4).Value = "-"
$ws.Cells.Item(8, 5).Value = "-"
$ws.Cells.Item(8, 6).Value = "-"

# --- Row 9 (new): the lunch break, now starting at 12:20 ---
$ws.Cells.Item(9, 1).Value = "12:20"
$ws.Cells.Item(9, 2).Value = "Almoço"
$ws.Cells.Item(9, 3).Value = "Almoço"
$ws.Cells.Item(9, 4).Value = "Almoço"
$ws.Cells.Item(9, 5).Value = "Almoço"
$ws.Cells.Item(9, 6).Value = "Almoço"

# --- Re-stamp the time column for the rest of the shifted rows ---
$ws.Cells.Item(10, 1).Value = "13:00"
$ws.Cells.Item(11, 1).Value = "13:50"
$ws.Cells.Item(12, 1).Value = "14:40"
$ws.Cells.Item(13, 1).Value = "15:30"
$ws.Cells.Item(14, 1).Value = "15:50"

# --- Append two new trailing rows for the extended schedule ---
$ws.Cells.Item(15, 1).Value = "16:40"
$ws.Cells.Item(15, 2).Value = "-"
$ws.Cells.Item(15, 3).Value = "-"
$ws.Cells.Item(15, 4).Value = "-"
$ws.Cells.Item(15, 5).Value = "-"
$ws.Cells.Item(15, 6).Value = "-"

$ws.Cells.Item(16, 1).Value = "17:30"
$ws.Cells.Item(16, 2).Value = "-"
$ws.Cells.Item(16, 3).Value = "-"
$ws.Cells.Item(16, 4).Value = "-"
$ws.Cells.Item(16, 5).Value = "-"
$ws.Cells.Item(16, 6).Value = "-"

$ws.Cells.Item(17, 1).Value = "18:20"
# Materialize the (intentionally blank) trailing cells B17:F17 so the row's
# full width is preserved, matching the rest of the table's shape.
$ws.Range("B17:F17").Borders.LineStyle = 0
